$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before the old "BARCODE" column (L) so that
# L:M become new tier-3 columns (JML3 / HARGAJUAL3), shifting the old
# L,M,N (BARCODE, NAMASUPPLIER, RAK) to N,O,P.
$ws.Range("L1:M1").EntireColumn.Insert()

# Header labels for the two new columns.
$ws.Range("L1").Value = "JML3"
$ws.Range("M1").Value = "HARGAJUAL3"

# Give the new header cells their own (new) fill style - a light gray
# accent, matching the other tier-header cells' shaded look.
$ws.Range("L1:M1").Interior.ThemeColor = 7
$ws.Range("L1:M1").Interior.TintAndShade = 0.6

# Update the tiered pricing data on row 2:
#  - HARGAJUAL1 (I2) goes up to 33000
#  - HARGAJUAL2 (K2) becomes 31000
#  - new JML3 (L2) = 100
#  - new HARGAJUAL3 (M2) = 28000 (the old HARGAJUAL2 value)
$ws.Range("I2").Value = 33000
$ws.Range("K2").Value = 31000
$ws.Range("L2").Value = 100
$ws.Range("M2").Value = 28000

# Match the author's final selection.
$ws.Range("A2").Select() | Out-Null
